$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) values
# D-column values are forced to Text format to preserve exact string
# representation (avoids Excel auto-converting numeric-looking strings
# to numbers and losing formatting such as trailing zeros).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.255.59'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.786.31'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').Value = '  -1.32%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.57'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').Value = '  -2.99%  '

$ws.Range('E6').Value = '  +0.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3790'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').Value = '  -1.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.63'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').Value = '  -3.31%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3422'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').Value = '  -2.73%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.196'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').Value = '  -2.92%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07486'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').Value = '  -3.26%  '

$ws.Range('E12').Value = '  -0.17%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.86'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').Value = '  -3.03%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.466'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').Value = '  -2.19%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.788.14'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').Value = '  -1.35%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.081'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').Value = '  -1.87%  '

$ws.Range('E17').Value = '  -2.21%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06647'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').Value = '  -1.92%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.86'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').Value = '  -3.54%  '

$ws.Range('E20').Value = '  +0.14%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.634'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').Value = '  +1.51%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.34'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').Value = '  -2.50%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.253.75'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').Value = '  -0.49%  '

$ws.Range('E24').Value = '  -5.99%  '

$ws.Range('E25').Value = '  -2.54%  '

$ws.Range('E26').Value = '  +1.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.541'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').Value = '  -6.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.30'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').Value = '  -4.67%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '152.90'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').Value = '  -1.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.991.73'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').Value = '  -1.22%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '134.28'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').Value = '  -1.72%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.017'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').Value = '  -2.10%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.080'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').Value = '  -4.41%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08705'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').Value = '  -0.88%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.29'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').Value = '  -3.75%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.663'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').Value = '  -3.59%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6937'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').Value = '  -1.82%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.449'
$ws.Range('D38').ClearFormats()

$ws.Range('E38').Value = '  -3.21%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2205'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').Value = '  -2.55%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.819'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').Value = '  -2.02%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06334'
$ws.Range('D41').ClearFormats()

$ws.Range('E41').Value = '  -3.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.02339'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').Value = '  -2.95%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.237'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').Value = '  -1.70%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.46'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').Value = '  -3.36%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6509'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').Value = '  -1.34%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.839'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').Value = '  -5.10%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.150'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').Value = '  -1.49%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '129.40'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').Value = '  -2.37%  '

$ws.Range('E50').Value = '  -3.16%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.07'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').Value = '  -1.91%  '
